$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.058999066880645
$ws.Range("D2").Value = 1.065704068918184
$ws.Range("E2").Value = 1.065335469089387
$ws.Range("F2").Value = 1.076931002269229
$ws.Range("I2").Value = 1.054039457981204
$ws.Range("J2").Value = 1.06398737435064
$ws.Range("K2").Value = 1.06841659283319
$ws.Range("L2").Value = 1.068048985894156
$ws.Range("M2").Value = 1.079613636321092
$ws.Range("N2").Value = 1.024993864333589
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.060068733111796
$ws.Range("D3").Value = 1.066584573063032
$ws.Range("E3").Value = 1.066289121326675
$ws.Range("F3").Value = 1.077967431960021
$ws.Range("I3").Value = 1.054393713301809
$ws.Range("J3").Value = 1.064709550031878
$ws.Range("K3").Value = 1.069112384800285
$ws.Range("L3").Value = 1.068817671420381
$ws.Range("M3").Value = 1.080467128583135
$ws.Range("N3").Value = 1.025241999958958
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.06076106481599
$ws.Range("D4").Value = 1.067154460084713
$ws.Range("E4").Value = 1.066906676658377
$ws.Range("F4").Value = 1.078638676230948
$ws.Range("I4").Value = 1.054621727879948
$ws.Range("J4").Value = 1.065176425676544
$ws.Range("K4").Value = 1.069562110045205
$ws.Range("L4").Value = 1.069314914694158
$ws.Range("M4").Value = 1.081019388133394
$ws.Range("N4").Value = 1.025402249670515
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.061052165666161
$ws.Range("D5").Value = 1.067394074038763
$ws.Range("E5").Value = 1.067166410985789
$ws.Range("F5").Value = 1.078921011665952
$ws.Range("I5").Value = 1.054717294803244
$ws.Range("J5").Value = 1.065372599301814
$ws.Range("K5").Value = 1.069751054812806
$ws.Range("L5").Value = 1.069523920058818
$ws.Range("M5").Value = 1.081251555966278
$ws.Range("N5").Value = 1.02546954408151
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.061101045396643
$ws.Range("D6").Value = 1.067434308221991
$ws.Range("E6").Value = 1.067210028196548
$ws.Range("F6").Value = 1.078968425478181
$ws.Range("I6").Value = 1.054733323875559
$ws.Range("J6").Value = 1.065405531792273
$ws.Range("K6").Value = 1.069782772440115
$ws.Range("L6").Value = 1.069559010874978
$ws.Range("M6").Value = 1.081290537833267
$ws.Range("N6").Value = 1.025480838730679
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060764954344474
$ws.Range("D7").Value = 1.067157661687325
$ws.Range("E7").Value = 1.066910146793595
$ws.Range("F7").Value = 1.078642448244277
$ws.Range("I7").Value = 1.054623005991166
$ws.Range("J7").Value = 1.065179047354912
$ws.Range("K7").Value = 1.069564635205296
$ws.Range("L7").Value = 1.069317707574818
$ws.Range("M7").Value = 1.081022490380019
$ws.Range("N7").Value = 1.025403149154978
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.059360527077663
$ws.Range("D8").Value = 1.066001609699896
$ws.Range("E8").Value = 1.06565766093481
$ws.Range("F8").Value = 1.077281142983803
$ws.Range("I8").Value = 1.054159431133415
$ws.Range("J8").Value = 1.064231523867485
$ws.Range("K8").Value = 1.068651842037631
$ws.Range("L8").Value = 1.068308797126121
$ws.Range("M8").Value = 1.079902079237018
$ws.Range("N8").Value = 1.025077787126338
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.056887170031495
$ws.Range("D9").Value = 1.063965609022504
$ws.Range("E9").Value = 1.063454308018611
$ws.Range("F9").Value = 1.0748870035212
$ws.Range("I9").Value = 1.053333275911329
$ws.Range("J9").Value = 1.062558659094076
$ws.Range("K9").Value = 1.06703957925276
$ws.Range("L9").Value = 1.066529852605119
$ws.Range("M9").Value = 1.077927735432594
$ws.Range("N9").Value = 1.024502087452323
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.055239216141092
$ws.Range("D10").Value = 1.062609051428427
$ws.Range("E10").Value = 1.061987912623595
$ws.Range("F10").Value = 1.073294065975983
$ws.Range("I10").Value = 1.052776277491131
$ws.Range("J10").Value = 1.061441273350205
$ws.Range("K10").Value = 1.065962197361694
$ws.Range("L10").Value = 1.065343158633701
$ws.Range("M10").Value = 1.076611503385702
$ws.Range("N10").Value = 1.024116705374987
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.054525857144981
$ws.Range("D11").Value = 1.062021836664393
$ws.Range("E11").Value = 1.061353546857765
$ws.Range("F11").Value = 1.072605059612632
$ws.Range("I11").Value = 1.052533615258268
$ws.Range("J11").Value = 1.0609569270653
$ws.Range("K11").Value = 1.065495079808368
$ws.Range("L11").Value = 1.064829137034378
$ws.Range("M11").Value = 1.07604156355181
$ws.Range("N11").Value = 1.023949457553781
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.054260915982127
$ws.Range("D12").Value = 1.061803746909199
$ws.Range("E12").Value = 1.061118004607151
$ws.Range("F12").Value = 1.072349244360364
$ws.Range("I12").Value = 1.052443257664677
$ws.Range("J12").Value = 1.060776942630721
$ws.Range("K12").Value = 1.065321480965155
$ws.Range("L12").Value = 1.064638180381657
$ws.Range("M12").Value = 1.075829862343568
$ws.Range("N12").Value = 1.023887278096162
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.054317745311844
$ws.Range("D13").Value = 1.061850526631152
$ws.Range("E13").Value = 1.061168525166803
$ws.Range("F13").Value = 1.072404112509761
$ws.Range("I13").Value = 1.052462649721934
$ws.Range("J13").Value = 1.060815553370325
$ws.Range("K13").Value = 1.065358722602734
$ws.Range("L13").Value = 1.064679142411885
$ws.Range("M13").Value = 1.075875272975782
$ws.Range("N13").Value = 1.023900618340374
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.054503956372119
$ws.Range("D14").Value = 1.062003808724821
$ws.Range("E14").Value = 1.061334075035383
$ws.Range("F14").Value = 1.072583911549659
$ws.Range("I14").Value = 1.052526150800538
$ws.Range("J14").Value = 1.060942051059045
$ws.Range("K14").Value = 1.065480731920751
$ws.Range("L14").Value = 1.064813353030117
$ws.Range("M14").Value = 1.07602406426621
$ws.Range("N14").Value = 1.023944318925428
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.05461869140059
$ws.Range("D15").Value = 1.062098254596815
$ws.Range("E15").Value = 1.061436087634841
$ws.Range("F15").Value = 1.07269470657377
$ws.Range("I15").Value = 1.052565246490614
$ws.Range("J15").Value = 1.061019980304884
$ws.Range("K15").Value = 1.065555893890105
$ws.Range("L15").Value = 1.064896041162912
$ws.Range("M15").Value = 1.076115739474214
$ws.Range("N15").Value = 1.023971236859451
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.055286563950925
$ws.Range("D16").Value = 1.062648026838395
$ws.Range("E16").Value = 1.062030025914735
$ws.Range("F16").Value = 1.073339808803866
$ws.Range("I16").Value = 1.05279235104398
$ws.Range("J16").Value = 1.061473407062351
$ws.Range("K16").Value = 1.065993185702507
$ws.Range("L16").Value = 1.065377268877583
$ws.Range("M16").Value = 1.07664932837821
$ws.Range("N16").Value = 1.024127797177019
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.055705560344305
$ws.Range("D17").Value = 1.06299293386103
$ws.Range("E17").Value = 1.06240274673353
$ws.Range("F17").Value = 1.073744664350989
$ws.Range("I17").Value = 1.052934412000013
$ws.Range("J17").Value = 1.061757692976723
$ws.Range("K17").Value = 1.066267325693216
$ws.Range("L17").Value = 1.065679083515465
$ws.Range("M17").Value = 1.076984034042821
$ws.Range("N17").Value = 1.024225903075328
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.055949974704631
$ws.Range("D18").Value = 1.063194129984643
$ws.Range("E18").Value = 1.06262020568613
$ws.Range("F18").Value = 1.073980881710499
$ws.Range("I18").Value = 1.053017131119628
$ws.Range("J18").Value = 1.061923462841388
$ws.Range("K18").Value = 1.066427168614698
$ws.Range("L18").Value = 1.065855109962765
$ws.Range("M18").Value = 1.07717926192131
$ws.Range("N18").Value = 1.024283090409739
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.056033317231716
$ws.Range("D19").Value = 1.063262735653274
$ws.Range("E19").Value = 1.062694363262225
$ws.Range("F19").Value = 1.074061437923834
$ws.Range("I19").Value = 1.053045312007224
$ws.Range("J19").Value = 1.06197997771644
$ws.Range("K19").Value = 1.066481661026837
$ws.Range("L19").Value = 1.0659151275559
$ws.Range("M19").Value = 1.077245829529655
$ws.Range("N19").Value = 1.024302583672393
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.055660603841823
$ws.Range("D20").Value = 1.062955926768906
$ws.Range("E20").Value = 1.062362751397619
$ws.Range("F20").Value = 1.073701219736035
$ws.Range("I20").Value = 1.052919184956602
$ws.Range("J20").Value = 1.061727196887066
$ws.Range("K20").Value = 1.066237919095267
$ws.Range("L20").Value = 1.065646703410316
$ws.Range("M20").Value = 1.076948123319157
$ws.Range("N20").Value = 1.024215380981711
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.054449120971543
$ws.Range("D21").Value = 1.061958670194081
$ws.Range("E21").Value = 1.061285322253662
$ws.Range("F21").Value = 1.072530962108912
$ws.Range("I21").Value = 1.052507457444445
$ws.Range("J21").Value = 1.060904802757967
$ws.Range("K21").Value = 1.065444805716813
$ws.Range("L21").Value = 1.064773832068325
$ws.Range("M21").Value = 1.07598024894279
$ws.Range("N21").Value = 1.023931451741739
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.053687599356303
$ws.Range("D22").Value = 1.061331817319192
$ws.Range("E22").Value = 1.060608417424814
$ws.Range("F22").Value = 1.071795825087621
$ws.Range("I22").Value = 1.052247303413965
$ws.Range("J22").Value = 1.060387287453033
$ws.Range("K22").Value = 1.064945619843801
$ws.Range("L22").Value = 1.06422487204633
$ws.Range("M22").Value = 1.075371706993879
$ws.Range("N22").Value = 1.023752609248862
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.054091278867785
$ws.Range("D23").Value = 1.061664108418495
$ws.Range("E23").Value = 1.06096720830518
$ws.Range("F23").Value = 1.072185473382754
$ws.Range("I23").Value = 1.052385337734792
$ws.Range("J23").Value = 1.060661674111661
$ws.Range("K23").Value = 1.065210297332685
$ws.Range("L23").Value = 1.064515900405485
$ws.Range("M23").Value = 1.075694306606273
$ws.Range("N23").Value = 1.023847447787282
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.055680917685542
$ws.Range("D24").Value = 1.062972648626164
$ws.Range("E24").Value = 1.062380823391523
$ws.Range("F24").Value = 1.073720850265035
$ws.Range("I24").Value = 1.052926065842726
$ws.Range("J24").Value = 1.061740976909529
$ws.Range("K24").Value = 1.066251206851139
$ws.Range("L24").Value = 1.065661334638266
$ws.Range("M24").Value = 1.076964349830187
$ws.Range("N24").Value = 1.024220135574663
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.05752642409519
$ws.Range("D25").Value = 1.064491829350048
$ws.Range("E25").Value = 1.064023487781688
$ws.Range("F25").Value = 1.075505391710808
$ws.Range("I25").Value = 1.053547955289997
$ws.Range("J25").Value = 1.062991513398818
$ws.Range("K25").Value = 1.067456836846794
$ws.Range("L25").Value = 1.066989882636213
$ws.Range("M25").Value = 1.078438152712101
$ws.Range("N25").Value = 1.024651199336978
